# Update Name of Algo
# Apply updated computed values to result_data_RandomForest sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.43729999999999
$ws.Range("B3").Value = 6.027099999999989
$ws.Range("D5").Value = -8.669299999999991
$ws.Range("E5").Value = 12.30269999999999
$ws.Range("E9").Value = 13.55620000000001
$ws.Range("E11").Value = 13.7274
$ws.Range("B14").Value = 9.093000000000004
$ws.Range("B21").Value = 5.724099999999996
$ws.Range("E21").Value = 13.06639999999999
$ws.Range("B23").Value = 5.665499999999998
$ws.Range("B25").Value = 5.836699999999992
